$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.240.14'
$ws.Range('E2').Value = '  +0.50%  '
$ws.Range('D3').Value = '2.642.90'
$ws.Range('E3').Value = '  +0.58%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'598.51"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('D6').Value = "'154.64"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.62%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = "'0.544"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.67%  '
$ws.Range('D9').Value = '2.641.92'
$ws.Range('E9').Value = '  +0.60%  '
$ws.Range('E10').Value = '  +8.89%  '
$ws.Range('E11').Value = '  -0.64%  '
$ws.Range('E12').Value = '  +0.79%  '
$ws.Range('E13').Value = '  +1.98%  '
$ws.Range('E14').Value = '  +2.66%  '
$ws.Range('D15').Value = "'27.88"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.15%  '
$ws.Range('D16').Value = '3.125.88'
$ws.Range('E16').Value = '  +0.65%  '
$ws.Range('D17').Value = '68.211.92'
$ws.Range('E17').Value = '  +0.73%  '
$ws.Range('D18').Value = '2.642.92'
$ws.Range('E18').Value = '  +0.88%  '
$ws.Range('E19').Value = '  -0.03%  '
$ws.Range('D20').Value = "'363.85"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.07%  '
$ws.Range('E21').Value = '  +0.38%  '
$ws.Range('D22').Value = "'4.35"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.77%  '
$ws.Range('D23').Value = "'4.84"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.03%  '
$ws.Range('E24').Value = '  -1.34%  '
$ws.Range('D25').Value = "'75.38"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.67%  '
$ws.Range('D26').Value = "'0.999"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('D27').Value = "'9.78"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.42%  '
$ws.Range('D30').Value = "'1.00"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('D31').Value = "'563.00"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.16%  '
$ws.Range('D32').Value = "'8.04"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.59%  '
$ws.Range('E33').Value = '  +0.18%  '
$ws.Range('E34').Value = '  +0.90%  '
$ws.Range('E35').Value = '  +2.19%  '
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('E37').Value = '  +2.98%  '
$ws.Range('D38').Value = "'161.84"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.41%  '
$ws.Range('D39').Value = "'19.31"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.90%  '
$ws.Range('D40').Value = "'0.373"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.53%  '
$ws.Range('D41').Value = "'1.88"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('D42').Value = "'5.33"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.53%  '
$ws.Range('D43').Value = '0.0₆0338'
$ws.Range('E43').Value = '  +0.64%  '
$ws.Range('E44').Value = '  -0.94%  '
$ws.Range('D45').Value = "'17.74"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.16%  '
$ws.Range('D46').Value = "'40.63"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.19%  '
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').Value = "'157.30"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.42%  '
$ws.Range('E49').Value = '  +1.89%  '
$ws.Range('E50').Value = '  +0.20%  '
$ws.Range('D51').Value = "'21.78"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.81%  '
